$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the 18 rows (film_id 1010101 .. 1060106b) that were removed from
#    the filtered export. They are currently contiguous at rows 151-168.
$ws.Range("A151:A168").EntireRow.Delete()

# 2) The film_id column (A) was exported as text (inline strings) but should
#    be numeric. Convert every remaining data row's A cell (rows 2..179)
#    from text to a true number, preserving the same digits.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) {
        $text = $cell.Value2
    }
    if ($text -ne $null -and $text -ne "") {
        $cell.Value2 = [double]$text
    }
}
